# Apply the edits described by the diff to the active document.
$d = $word.ActiveDocument

# 1) "Due to a lot of noise ... standardized per 100K ..." ->
#    "Due to noise ... standardize per 100K ..."
$d.Content.Find.Execute(
    "Due to a lot of noise in the data regarding state size in acres or state size in population it appears that new scaling columns would be useful to standardized per 100K of state population",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Due to noise in the data regarding state size in acres or state size in population it appears that new scaling columns would be useful to standardize per 100K of state population",
    2) | Out-Null

# 4) "Snow Making_ac" -> "Snow_Making_ac"
$d.Content.Find.Execute(
    "Snow Making_ac",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Snow_Making_ac",
    2) | Out-Null

# 5) Insert "and " between "linear regression model" and "the random forest model"
$d.Content.Find.Execute(
    "When comparing the linear regression model the random forest model the random forest model has a lower",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "When comparing the linear regression model and the random forest model the random forest model has a lower",
    2) | Out-Null

# 6) Append new sentence after "...actual ticket price is $81."
$d.Content.Find.Execute(
    "In the final analysis, the predicted ticket price generated by the model is `$95.87 and the actual ticket price is `$81.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In the final analysis, the predicted ticket price generated by the model is `$95.87 and the actual ticket price is `$81.  It is reasonable to assume that some of the other resorts used in deriving the predictive model under charged and some over charged thus evening out the predicted price of `$95.87.",
    2) | Out-Null
